$d = $word.ActiveDocument

# Change 1: "817 " (home address) -> "39 "
# The "817" run is followed by a separate space run with special spacing.
# Target just the home-address paragraph (second paragraph under "HOME ADDRESS").
$found = $d.Content.Find.Execute("817 Birch", $true, $false, $false, $false, $false,
                         $true, 1, $false, "39 Birch", 2)
Write-Output "Found: $found"
